# NIT-9016953574.xlsx — "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker table (rows 16-18) held 3 account-statement lines (one per
# employee, all for period 2506). The update removes those and adds new
# lines covering two periods (2506 and 2507) for the same three employees,
# ordered JULIO CESAR RINCON ACOSTA, KATIA RUIZ HERRERA, ANDRES ALONSO
# PAYARES CARDOZA (each appearing twice, once per period) -> 6 data rows.
# The totals (Valor Mora / Cant. Periodos) and the signature block (which
# shifts down 3 rows) are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Grow the detail table from 3 to 6 rows -----------------------
# Insert 3 blank rows right after the current last detail row (18) so the
# trailing signature block (currently rows 23/24) is pushed down to 26/27,
# exactly like the authored change.
$ws.Range("B19:J21").EntireRow.Insert()

# Re-apply the correct formatting to the (now 6) detail rows:
#  - rows 16-20 use the "normal" row style (thin borders all round)
#  - row 21 (the new last row) uses the "closing" style with the heavier
#    bottom border, same as row 18 used to have before the insert.
$ws.Range("B18:J18").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B17:J17").Copy()
$ws.Range("B18:J20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Write the new detail rows -------------------------------------
# Columns: B Tipo Doc | C N Doc | D Nombre | E Periodo Mora | F Valor Mora | G Salario Basico
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73194195"
$ws.Range("D16").Value = "JULIO CESAR RINCON ACOSTA"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73194195"
$ws.Range("D17").Value = "JULIO CESAR RINCON ACOSTA"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "45475098"
$ws.Range("D18").Value = "KATIA RUIZ HERRERA"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 128000
$ws.Range("G18").Value = 3200000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45475098"
$ws.Range("D19").Value = "KATIA RUIZ HERRERA"
$ws.Range("E19").Value = "2506"
$ws.Range("F19").Value = 128000
$ws.Range("G19").Value = 3200000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1050968370"
$ws.Range("D20").Value = "ANDRES ALONSO PAYARES CARDOZA"
$ws.Range("E20").Value = "2507"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1050968370"
$ws.Range("D21").Value = "ANDRES ALONSO PAYARES CARDOZA"
$ws.Range("E21").Value = "2506"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

# --- 3. Update the summary block --------------------------------------
# Valor Mora total = sum of the 6 new detail lines; Cant. Periodos = 2
$ws.Range("E11").Value = 483760
$ws.Range("F13").Value = 2

Write-Host "edit applied"
